# Add row 7: new SmartScore submission from Omar Huerta (Streamlit export)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values must stay literal text (e.g. numeric-looking SmartScore
# strings and the blank Grupo_Experimental cell) are pre-formatted as Text so
# Excel does not auto-convert them to numbers on entry.
$textCols = @("B7", "I7", "L7", "O7", "R7", "U7", "X7", "AA7", "AD7", "AG7")
foreach ($addr in $textCols) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A7").Value = "Omar Huerta_20251119_152531"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "Omar Huerta"
$ws.Range("D7").Value = 42
$ws.Range("E7").Value = "Male"
$ws.Range("F7").Value = "2025-11-19 15:25:32"
$ws.Range("G7").Value = "{`n  `"portion`": 0.8,`n  `"diet`": 0.8571428571428571,`n  `"salt`": 0.6,`n  `"fat`": 0.4,`n  `"natural`": 1.0,`n  `"convenience`": 1.0,`n  `"price`": 0.8`n}"
$ws.Range("H7").Value = "Nongshim Neoguri Spicy Seafood"
$ws.Range("I7").Value = "0.602"
$ws.Range("J7").Value = "Sabor a marisco, umami, picante equilibrado, buena textura, algo salado"
$ws.Range("K7").Value = "Nissin Chow Mein Teriyaki Beef"
$ws.Range("L7").Value = "0.497"
$ws.Range("M7").Value = "Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa"
$ws.Range("N7").Value = "Maruchan Ramen Sabor Pollo"
$ws.Range("O7").Value = "0.449"
$ws.Range("P7").Value = "Sabor clásico, económico, alto en sodio, no saludable, nostálgico"
$ws.Range("Q7").Value = "Amy’s Macaroni & Cheese (frozen)"
$ws.Range("R7").Value = "0.643"
$ws.Range("S7").Value = "Queso real, textura casera, sin conservadores, alto en grasa, algo caro"
$ws.Range("T7").Value = "Kraft Macaroni & Cheese Dinner"
$ws.Range("U7").Value = "0.571"
$ws.Range("V7").Value = "Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato"
$ws.Range("W7").Value = "Annie’s Shells & White Cheddar"
$ws.Range("X7").Value = "0.522"
$ws.Range("Y7").Value = "Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños"
$ws.Range("Z7").Value = "Wild Planet Wild Tuna Pasta Salad"
$ws.Range("AA7").Value = "0.734"
$ws.Range("AB7").Value = "Sabor fresco, buena proteína, saludable, porción algo pequeña"
$ws.Range("AC7").Value = "StarKist Chicken Creations (Chicken Salad)"
$ws.Range("AD7").Value = "0.553"
$ws.Range("AE7").Value = "Portátil, saludable, fácil, buena textura, sabor suave"
$ws.Range("AF7").Value = "Kitchens of India Variety Pack"
$ws.Range("AG7").Value = "0.540"
$ws.Range("AH7").Value = "Sabor auténtico, variedad, vegetariano, necesita arroz o pan, buena calidad"

# Excel auto-grows row height to fit the multi-line JSON blob in column G;
# restore the sheet's default height so the new row matches the others.
$ws.Rows.Item(7).RowHeight = 15
